$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting of the last existing row (85) into the new row (86)
$ws.Range("A85:L85").Copy()
$ws.Range("A86:L86").PasteSpecial(-4104)
$ws.Rows.Item(86).RowHeight = $ws.Rows.Item(85).RowHeight

# Fill in the new training record values, matching the order in which the
# author appears to have entered them (so new shared strings land at the
# same indices as in the target workbook)
$ws.Range("A86").Value = 43223.375
$ws.Range("B86").Value = "分类O"
$ws.Range("C86").Value = "14分类"
$ws.Range("F86").Value = 0.67
$ws.Range("G86").Value = 0.61
$ws.Range("H86").Value = 0.999
$ws.Range("I86").Value = 0.98
$ws.Range("E86").Value = "最高标签，重新训练，原始数据加PCA及Wavelet处理数据train-hjxh365-2018-4-16-day-high-original-pca99-wavelet20"
$ws.Range("J86").Value = "经过约15小时，拟合精度和泛化精度开始接近，后面分化，泛化精度后来下降，最好达到0.64。"
$ws.Range("L86").Value = "logs-hjxh-2018-5-3-high-original-pca99-wavelet20-percent64"
$ws.Range("K86").Value = "python feed_run.py --output_mode=classes --output_nodes=14 --input_nums=96 --input_nodes=96 --low_nums=2 --low_nodes=96 --low_fun=elu --one_hot=True --input_fun=tanh --batch_size=100 --learning_rate=0.001 --train_mode=Adadelta --eval_size=5400 --test_size=1339 --use_biases=yes  --use_bn_input=True --dropout_low=0.8"
$ws.Range("D86").Value = "batch_size=100 low_nums=2 use_biases=yes use_bn_input=True   dropout_low=0.8 "

# Update the view: scroll window and selection as in the author's edit
$ws.Application.ActiveWindow.ScrollRow = 85
$ws.Range("E69").Select()
